# Rename the three sheets affected by the Product -> Beverage terminology change.
$wb = $excel.ActiveWorkbook

$wsProductTopping  = $wb.Worksheets.Item("ProductTopping")
$wsProductCategory = $wb.Worksheets.Item("ProductCategory")
$wsProduct         = $wb.Worksheets.Item("Product")

$wsProductTopping.Name  = "BeverageTopping"
$wsProductCategory.Name = "BeverageCategory"
$wsProduct.Name         = "Beverage"

# --- Update header labels that referenced the old "product..." names ---

# BeverageCategory sheet: A1/B1 headers
$wsBeverageCategory = $wb.Worksheets.Item("BeverageCategory")
$wsBeverageCategory.Range("A1").Value = "BeverageCategoryNo"
$wsBeverageCategory.Range("B1").Value = "BeverageCategoryName"

# Beverage sheet: A1/B1/C1 headers
$wsBeverage = $wb.Worksheets.Item("Beverage")
$wsBeverage.Range("A1").Value = "BeverageNo"
$wsBeverage.Range("B1").Value = "BeverageName"
$wsBeverage.Range("C1").Value = "BeverageCategoryNo"

# --- Restore per-sheet view state (selection / zoom) to match the target workbook ---

$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("U58").Select()

$wsBeverageTopping = $wb.Worksheets.Item("BeverageTopping")
$wsBeverageTopping.Activate()
$wsBeverageTopping.Range("F16").Select()

$wsToppingCategory = $wb.Worksheets.Item("ToppingCategory")
$wsToppingCategory.Activate()
$excel.ActiveWindow.Zoom = 142
$wsToppingCategory.Range("H19").Select()

$wsTopping = $wb.Worksheets.Item("Topping")
$wsTopping.Activate()
$wsTopping.Range("A46").Select()

$wsSubTopping = $wb.Worksheets.Item("SubTopping")
$wsSubTopping.Activate()
$wsSubTopping.Range("E28").Select()

$wsBeverageCategory.Activate()
$wsBeverageCategory.Range("B2").Select()

$wsBeverage.Activate()
$wsBeverage.Range("C15").Select()
